$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (dates 2021-10-25 .. 2021-10-31), matching the
# existing covid_totals data layout:
#   A: date (text)            B: areaType        C: areaCode
#   D: areaName                E: cumCases        F: newCases
#   G: newDeaths28Days         H: cumDeaths28Days
$rows = @(
    @{ Row=440; Date="2021-10-25"; E=8809774; F=36567; G=38;  H=139571 },
    @{ Row=441; Date="2021-10-26"; E=8853227; F=40954; G=263; H=139834 },
    @{ Row=442; Date="2021-10-27"; E=8897149; F=43941; G=207; H=140041 },
    @{ Row=443; Date="2021-10-28"; E=8936155; F=39842; G=165; H=140206 },
    @{ Row=444; Date="2021-10-29"; E=8979236; F=43467; G=186; H=140392 },
    @{ Row=445; Date="2021-10-30"; E=9019962; F=41278; G=166; H=140558 },
    @{ Row=446; Date="2021-10-31"; E=9057629; F=38009; G=74;  H=140632 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A holds a date-formatted string ("2021-10-25"). A plain
    # assignment gets auto-parsed into a date serial by Excel, so force
    # text entry with a leading apostrophe and then strip the resulting
    # "quote prefix" cell style so no stray formatting is left behind.
    $ws.Cells.Item($row, 1).Value = "'" + $r.Date
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = "overview"
    $ws.Cells.Item($row, 3).Value = "K02000001"
    $ws.Cells.Item($row, 4).Value = "United Kingdom"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
}

Write-Host "Added rows 440-446"
